$wb = $excel.ActiveWorkbook

function Set-LiteralText($sheet, $range, $text) {
    # Writing a numeric/percent/date-looking string straight to .Value lets
    # Excel's smart-parsing reinterpret it (e.g. "20.0%" -> 0.2, "6.6" -> 6.6,
    # "9873392286" -> a number). Route it through a scratch cell holding a
    # text-literal formula, then paste only the resulting text value in, so
    # the destination keeps its existing style and becomes a real string.
    $scratch = $sheet.Range("ZZ100")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $scratch.Value = ""
    $excel.CutCopyMode = 0
}

function Copy-Format($sheet, $srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
}

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")

# Email row: 0/5 -> 1/5, 0.0% -> 20.0% (now has a mismatch, so it gets the
# same "bad" highlight style as the Name row, B2/C2 style group).
Set-LiteralText $summary $summary.Range("B3") "1/5"
Copy-Format $summary $summary.Range("B2") $summary.Range("B3")
Set-LiteralText $summary $summary.Range("C3") "20.0%"

# Mobile row: 0/5 -> 1/5, 0.0% -> 20.0%
Set-LiteralText $summary $summary.Range("B4") "1/5"
Copy-Format $summary $summary.Range("B2") $summary.Range("B4")
Set-LiteralText $summary $summary.Range("C4") "20.0%"

# Experience row: 4/5 -> 1/5, 80.0% -> 20.0% (style is unchanged)
Set-LiteralText $summary $summary.Range("B6") "1/5"
Set-LiteralText $summary $summary.Range("C6") "20.0%"

# Overall wrong cells: 8/30 -> 7/30, 26.7% -> 23.3%
Set-LiteralText $summary $summary.Range("B9") "7/30"
Set-LiteralText $summary $summary.Range("C9") "23.3%"

# ---- Report sheet ----
$report = $wb.Worksheets.Item("Report")

# Row 2: extracted experience 7.7 -> 6.6, which now matches the expected
# value, so it gets the "match" style (same as K4).
Set-LiteralText $report $report.Range("K2") "6.6"
Copy-Format $report $report.Range("K4") $report.Range("K2")

# Row 4: extracted email/mobile are now wrong, so they get the "mismatch"
# style (same as C2).
Set-LiteralText $report $report.Range("E4") "malikvibhor@linkedin.com"
Copy-Format $report $report.Range("C2") $report.Range("E4")

Set-LiteralText $report $report.Range("G4") "9873392286"
Copy-Format $report $report.Range("C2") $report.Range("G4")

# Rows 5 & 6: extracted experience now matches expected, so they get the
# "match" style (same as K4).
Copy-Format $report $report.Range("K4") $report.Range("K5")
Copy-Format $report $report.Range("K4") $report.Range("K6")
